$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the old generic placeholder headers (text1..text11) with the real
# field names, and add a new "url_imagen " header for column L.
$ws.Range("A1").Value = "puesto"
$ws.Range("B1").Value = "curp"
$ws.Range("C1").Value = "telefono"
$ws.Range("E1").Value = "alergia"
$ws.Range("F1").Value = "fecha_expedicion"
$ws.Range("D1").Value = "tipo_sangre"
$ws.Range("G1").Value = "fecha_vigencia"
$ws.Range("H1").Value = "familiar"
$ws.Range("I1").Value = "parentesco"
$ws.Range("J1").Value = "telefono_parentesco"
$ws.Range("K1").Value = "nombre_elemento"

# New column L gets the same header style as the rest of row 1, plus its text.
$ws.Range("K1").Copy()
$ws.Range("L1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Range("L1").Value = "url_imagen "

# Let column widths auto-fit the new header text. Only the columns whose new
# header is now longer than their existing data actually need to grow; the
# rest keep their original best-fit width. The brand-new column L is left at
# Excel's default width (never explicitly auto-fit).
$ws.Columns("D:D").AutoFit()
$ws.Columns("F:G").AutoFit()
$ws.Columns("I:J").AutoFit()

# Scroll/select as in the edited workbook: viewport starts at column C, with
# the freshly added L1 header area in focus.
$ws.Range("L2").Select()
$ws.Application.ActiveWindow.ScrollColumn = 3
